# New ids for series and indicators
#
# The catalog's INDICATOR_LABEL column (F) for the set of rows belonging
# to the I / IV / V domain series (rows 89-123) is re-pointed to new
# "provisional" labels that are the original label with a "q" suffix
# (e.g. "I.1" -> "I.1q", "V.4" -> "V.4q"). Writing the values through the
# object model lets Excel create/dedupe the new shared-string entries the
# same way the original edit did.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalog_series_old")

for ($r = 89; $r -le 123; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = $cell.Value() + "q"
}

# Leave the sheet with the same selection state captured after the edit.
[void]$ws.Range("F126:F127").Select()
